# Daily attendance processing - 2025-10-13 05:45:45
# Normalizes the "Recorded By" column (G) so that the "System" entry is
# always listed first among the recorders. When a cell has no "System"
# entry, the existing entries are reversed (matches the two-recorder
# admin/user swap observed in the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    if ($hasSystem) {
        if ($parts[0].Equals("System")) { continue }

        $rest = @()
        $removed = $false
        foreach ($p in $parts) {
            if ((-not $removed) -and $p.Equals("System")) {
                $removed = $true
                continue
            }
            $rest += $p
        }
        $newParts = @("System") + $rest
    } else {
        $newParts = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $cell.Value = [string]::Join(", ", $newParts)
}
